$d = $word.ActiveDocument

# The new "DICAS DE PREPARO:" section is appended at the very end of the
# document body, right after the "OBSERVAÇÃO: ... esfriar." paragraph.
$rng = $d.Content
$rng.Collapse(0)  # wdCollapseEnd

# Build the three new paragraphs ("DICAS DE PREPARO:" section) as raw
# WordprocessingML and insert them using InsertXML so that the exact
# run/paragraph formatting (including complex-script bold/font flags)
# matches what Word itself would produce.
$body = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">DICAS DE PREPARO: </w:t></w:r></w:p>' + `
        '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Para um bolo mais fofo, peneire a farinha de trigo.</w:t></w:r></w:p>' + `
        '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Você pode adicionar nozes picadas à massa para um toque especial.</w:t></w:r></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
          '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' + `
            '<pkg:xmlData>' + `
              '<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' + `
                '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' + `
              '</Relationships>' + `
            '</pkg:xmlData>' + `
          '</pkg:part>' + `
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
            '<pkg:xmlData>' + `
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
                '<w:body>' + $body + '</w:body>' + `
              '</w:document>' + `
            '</pkg:xmlData>' + `
          '</pkg:part>' + `
        '</pkg:package>'

$rng.InsertXML($xml) | Out-Null

Write-Output "done"
